$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(55, 1).Value = "'2026-02-01"
$ws.Cells.Item(55, 2).Value = "'18:29:48"
$ws.Cells.Item(55, 3).Value = "'18:00"
$ws.Cells.Item(55, 4).Value = "'Bathroom"
$ws.Cells.Item(55, 5).Value = "'Motion Detected"
$ws.Cells.Item(55, 6).Value = "'Active"
$ws.Cells.Item(56, 1).Value = "'2026-02-01"
$ws.Cells.Item(56, 2).Value = "'18:29:50"
$ws.Cells.Item(56, 3).Value = "'18:00"
$ws.Cells.Item(56, 4).Value = "'Bathroom"
$ws.Cells.Item(56, 5).Value = "'No Motion"
$ws.Cells.Item(56, 6).Value = "'Inactive"
$ws.Cells.Item(57, 1).Value = "'2026-02-01"
$ws.Cells.Item(57, 2).Value = "'18:29:53"
$ws.Cells.Item(57, 3).Value = "'18:00"
$ws.Cells.Item(57, 4).Value = "'Bathroom"
$ws.Cells.Item(57, 5).Value = "'Motion Detected"
$ws.Cells.Item(57, 6).Value = "'Active"
$ws.Cells.Item(58, 1).Value = "'2026-02-01"
$ws.Cells.Item(58, 2).Value = "'18:30:00"
$ws.Cells.Item(58, 3).Value = "'18:00"
$ws.Cells.Item(58, 4).Value = "'Bathroom"
$ws.Cells.Item(58, 5).Value = "'No Motion"
$ws.Cells.Item(58, 6).Value = "'Inactive"
$ws.Cells.Item(59, 1).Value = "'2026-02-01"
$ws.Cells.Item(59, 2).Value = "'18:30:01"
$ws.Cells.Item(59, 3).Value = "'18:00"
$ws.Cells.Item(59, 4).Value = "'Bathroom"
$ws.Cells.Item(59, 5).Value = "'Motion Detected"
$ws.Cells.Item(59, 6).Value = "'Active"
$ws.Cells.Item(60, 1).Value = "'2026-02-01"
$ws.Cells.Item(60, 2).Value = "'18:30:09"
$ws.Cells.Item(60, 3).Value = "'18:00"
$ws.Cells.Item(60, 4).Value = "'Bathroom"
$ws.Cells.Item(60, 5).Value = "'No Motion"
$ws.Cells.Item(60, 6).Value = "'Inactive"
$ws.Cells.Item(61, 1).Value = "'2026-02-01"
$ws.Cells.Item(61, 2).Value = "'18:30:14"
$ws.Cells.Item(61, 3).Value = "'18:00"
$ws.Cells.Item(61, 4).Value = "'Bathroom"
$ws.Cells.Item(61, 5).Value = "'No Motion"
$ws.Cells.Item(61, 6).Value = "'Inactive"
$ws.Cells.Item(62, 1).Value = "'2026-02-01"
$ws.Cells.Item(62, 2).Value = "'18:30:19"
$ws.Cells.Item(62, 3).Value = "'18:00"
$ws.Cells.Item(62, 4).Value = "'Bathroom"
$ws.Cells.Item(62, 5).Value = "'No Motion"
$ws.Cells.Item(62, 6).Value = "'Inactive"
$ws.Cells.Item(63, 1).Value = "'2026-02-01"
$ws.Cells.Item(63, 2).Value = "'18:30:24"
$ws.Cells.Item(63, 3).Value = "'18:00"
$ws.Cells.Item(63, 4).Value = "'Bathroom"
$ws.Cells.Item(63, 5).Value = "'Motion Detected"
$ws.Cells.Item(63, 6).Value = "'Active"
$ws.Cells.Item(64, 1).Value = "'2026-02-01"
$ws.Cells.Item(64, 2).Value = "'18:30:32"
$ws.Cells.Item(64, 3).Value = "'18:00"
$ws.Cells.Item(64, 4).Value = "'Bathroom"
$ws.Cells.Item(64, 5).Value = "'No Motion"
$ws.Cells.Item(64, 6).Value = "'Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Cells.Item(113, 1).Value = "'2026-02-01"
$ws.Cells.Item(113, 2).Value = "'18:29:48"
$ws.Cells.Item(113, 3).Value = "'18:00"
$ws.Cells.Item(113, 4).Value = "'Bathroom"
$ws.Cells.Item(113, 5).Value = "'79.4%"
$ws.Cells.Item(113, 6).Value = "'Active"
$ws.Cells.Item(114, 1).Value = "'2026-02-01"
$ws.Cells.Item(114, 2).Value = "'18:29:49"
$ws.Cells.Item(114, 3).Value = "'18:00"
$ws.Cells.Item(114, 4).Value = "'Bathroom"
$ws.Cells.Item(114, 5).Value = "'78.3%"
$ws.Cells.Item(114, 6).Value = "'Active"
$ws.Cells.Item(115, 1).Value = "'2026-02-01"
$ws.Cells.Item(115, 2).Value = "'18:29:59"
$ws.Cells.Item(115, 3).Value = "'18:00"
$ws.Cells.Item(115, 4).Value = "'Bathroom"
$ws.Cells.Item(115, 5).Value = "'78.2%"
$ws.Cells.Item(115, 6).Value = "'Active"
$ws.Cells.Item(116, 1).Value = "'2026-02-01"
$ws.Cells.Item(116, 2).Value = "'18:30:05"
$ws.Cells.Item(116, 3).Value = "'18:00"
$ws.Cells.Item(116, 4).Value = "'Bathroom"
$ws.Cells.Item(116, 5).Value = "'79.2%"
$ws.Cells.Item(116, 6).Value = "'Active"
$ws.Cells.Item(117, 1).Value = "'2026-02-01"
$ws.Cells.Item(117, 2).Value = "'18:30:09"
$ws.Cells.Item(117, 3).Value = "'18:00"
$ws.Cells.Item(117, 4).Value = "'Bathroom"
$ws.Cells.Item(117, 5).Value = "'78.3%"
$ws.Cells.Item(117, 6).Value = "'Active"
$ws.Cells.Item(118, 1).Value = "'2026-02-01"
$ws.Cells.Item(118, 2).Value = "'18:30:14"
$ws.Cells.Item(118, 3).Value = "'18:00"
$ws.Cells.Item(118, 4).Value = "'Bathroom"
$ws.Cells.Item(118, 5).Value = "'79.3%"
$ws.Cells.Item(118, 6).Value = "'Active"
$ws.Cells.Item(119, 1).Value = "'2026-02-01"
$ws.Cells.Item(119, 2).Value = "'18:30:19"
$ws.Cells.Item(119, 3).Value = "'18:00"
$ws.Cells.Item(119, 4).Value = "'Bathroom"
$ws.Cells.Item(119, 5).Value = "'78.5%"
$ws.Cells.Item(119, 6).Value = "'Active"
$ws.Cells.Item(120, 1).Value = "'2026-02-01"
$ws.Cells.Item(120, 2).Value = "'18:30:25"
$ws.Cells.Item(120, 3).Value = "'18:00"
$ws.Cells.Item(120, 4).Value = "'Bathroom"
$ws.Cells.Item(120, 5).Value = "'79.5%"
$ws.Cells.Item(120, 6).Value = "'Active"
$ws.Cells.Item(121, 1).Value = "'2026-02-01"
$ws.Cells.Item(121, 2).Value = "'18:30:29"
$ws.Cells.Item(121, 3).Value = "'18:00"
$ws.Cells.Item(121, 4).Value = "'Bathroom"
$ws.Cells.Item(121, 5).Value = "'78.6%"
$ws.Cells.Item(121, 6).Value = "'Active"
$ws.Cells.Item(122, 1).Value = "'2026-02-01"
$ws.Cells.Item(122, 2).Value = "'18:30:34"
$ws.Cells.Item(122, 3).Value = "'18:00"
$ws.Cells.Item(122, 4).Value = "'Bathroom"
$ws.Cells.Item(122, 5).Value = "'79.5%"
$ws.Cells.Item(122, 6).Value = "'Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Cells.Item(113, 1).Value = "'2026-02-01"
$ws.Cells.Item(113, 2).Value = "'18:29:48"
$ws.Cells.Item(113, 3).Value = "'18:00"
$ws.Cells.Item(113, 4).Value = "'Bathroom"
$ws.Cells.Item(113, 5).Value = "'29.6C"
$ws.Cells.Item(113, 6).Value = "'Active"
$ws.Cells.Item(114, 1).Value = "'2026-02-01"
$ws.Cells.Item(114, 2).Value = "'18:29:50"
$ws.Cells.Item(114, 3).Value = "'18:00"
$ws.Cells.Item(114, 4).Value = "'Bathroom"
$ws.Cells.Item(114, 5).Value = "'29.5C"
$ws.Cells.Item(114, 6).Value = "'Active"
$ws.Cells.Item(115, 1).Value = "'2026-02-01"
$ws.Cells.Item(115, 2).Value = "'18:30:00"
$ws.Cells.Item(115, 3).Value = "'18:00"
$ws.Cells.Item(115, 4).Value = "'Bathroom"
$ws.Cells.Item(115, 5).Value = "'29.5C"
$ws.Cells.Item(115, 6).Value = "'Active"
$ws.Cells.Item(116, 1).Value = "'2026-02-01"
$ws.Cells.Item(116, 2).Value = "'18:30:05"
$ws.Cells.Item(116, 3).Value = "'18:00"
$ws.Cells.Item(116, 4).Value = "'Bathroom"
$ws.Cells.Item(116, 5).Value = "'29.5C"
$ws.Cells.Item(116, 6).Value = "'Active"
$ws.Cells.Item(117, 1).Value = "'2026-02-01"
$ws.Cells.Item(117, 2).Value = "'18:30:10"
$ws.Cells.Item(117, 3).Value = "'18:00"
$ws.Cells.Item(117, 4).Value = "'Bathroom"
$ws.Cells.Item(117, 5).Value = "'29.5C"
$ws.Cells.Item(117, 6).Value = "'Active"
$ws.Cells.Item(118, 1).Value = "'2026-02-01"
$ws.Cells.Item(118, 2).Value = "'18:30:15"
$ws.Cells.Item(118, 3).Value = "'18:00"
$ws.Cells.Item(118, 4).Value = "'Bathroom"
$ws.Cells.Item(118, 5).Value = "'29.5C"
$ws.Cells.Item(118, 6).Value = "'Active"
$ws.Cells.Item(119, 1).Value = "'2026-02-01"
$ws.Cells.Item(119, 2).Value = "'18:30:20"
$ws.Cells.Item(119, 3).Value = "'18:00"
$ws.Cells.Item(119, 4).Value = "'Bathroom"
$ws.Cells.Item(119, 5).Value = "'29.5C"
$ws.Cells.Item(119, 6).Value = "'Active"
$ws.Cells.Item(120, 1).Value = "'2026-02-01"
$ws.Cells.Item(120, 2).Value = "'18:30:25"
$ws.Cells.Item(120, 3).Value = "'18:00"
$ws.Cells.Item(120, 4).Value = "'Bathroom"
$ws.Cells.Item(120, 5).Value = "'29.5C"
$ws.Cells.Item(120, 6).Value = "'Active"
$ws.Cells.Item(121, 1).Value = "'2026-02-01"
$ws.Cells.Item(121, 2).Value = "'18:30:30"
$ws.Cells.Item(121, 3).Value = "'18:00"
$ws.Cells.Item(121, 4).Value = "'Bathroom"
$ws.Cells.Item(121, 5).Value = "'29.4C"
$ws.Cells.Item(121, 6).Value = "'Active"
$ws.Cells.Item(122, 1).Value = "'2026-02-01"
$ws.Cells.Item(122, 2).Value = "'18:30:35"
$ws.Cells.Item(122, 3).Value = "'18:00"
$ws.Cells.Item(122, 4).Value = "'Bathroom"
$ws.Cells.Item(122, 5).Value = "'29.4C"
$ws.Cells.Item(122, 6).Value = "'Active"

$ws = $wb.Worksheets.Item("Proximity")
$ws.Cells.Item(44, 1).Value = "'2026-02-01"
$ws.Cells.Item(44, 2).Value = "'18:30:23"
$ws.Cells.Item(44, 3).Value = "'18:00"
$ws.Cells.Item(44, 4).Value = "'Bathroom Door"
$ws.Cells.Item(44, 5).Value = "'ENTER"
$ws.Cells.Item(44, 6).Value = "'User ENTERED Bathroom"
$ws.Cells.Item(45, 1).Value = "'2026-02-01"
$ws.Cells.Item(45, 2).Value = "'18:30:26"
$ws.Cells.Item(45, 3).Value = "'18:00"
$ws.Cells.Item(45, 4).Value = "'Bathroom Door"
$ws.Cells.Item(45, 5).Value = "'EXIT"
$ws.Cells.Item(45, 6).Value = "'User EXITED Bathroom"
$ws.Cells.Item(46, 1).Value = "'2026-02-01"
$ws.Cells.Item(46, 2).Value = "'18:30:30"
$ws.Cells.Item(46, 3).Value = "'18:00"
$ws.Cells.Item(46, 4).Value = "'Bathroom Door"
$ws.Cells.Item(46, 5).Value = "'ENTER"
$ws.Cells.Item(46, 6).Value = "'User ENTERED Bathroom"

